$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update parameter values (Min/Max columns for diffusion model)
$ws.Range("D2").Value = 0.1
$ws.Range("C3").Value = 5
$ws.Range("C4").Value = 0.01
$ws.Range("D4").Value = 0.1

# Update the active selection to C4
$ws.Range("C4").Select()
